$p = $ppt.ActivePresentation

# Remove the slide with SlideID="267" (the last slide, slide13.xml),
# which is an empty Title + Content Placeholder slide.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $s = $p.Slides.Item($i)
    if ($s.SlideID -eq 267) {
        $s.Delete()
        break
    }
}
